# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style used by the other header
# cells (e.g. AA1) by copying its format onto the new header cells.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-43 - same season record (91 wins, 72 losses, 0 ties) repeated
# for every player row.
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 91  # AD = column 30 = Wins
    $ws.Cells.Item($row, 31).Value = 72  # AE = column 31 = Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF = column 32 = Ties
}
